$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Remove all existing hyperlinks up front; they will be re-created below once
# the report rows have been regenerated, matching a freshly generated report.
$ws1.UsedRange.Hyperlinks.Delete()
$ws2.UsedRange.Hyperlinks.Delete()
$ws3.UsedRange.Hyperlinks.Delete()

# The localization request "f032c20d-7abb-4626-b4dd-30d58b520b37" has been
# handed off and is no longer pending, so its row is removed from every
# sheet. Deleting row 3 shifts the ".localization-config" row up from row 4
# to row 3 and keeps the original per-column styles intact.
$ws1.Rows(3).Delete()
$ws2.Rows(3).Delete()
$ws3.Rows(3).Delete()

# The remaining request (ebee5ea8-...) is now ready to be handed off again.
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws2.Range("B2").Value = "Ready for handoff"
$ws3.Range("B2").Value = "Ready for handoff"

# Refresh the handoff timestamps to reflect the new report generation time.
$ws2.Range("D2").Value = "2016-03-04 08:32:06"
$ws3.Range("D2").Value = "2016-03-04 08:32:17"

# Re-create the hyperlinks for the rows that remain, in report order, so
# that the relationship ids come out sequential starting right after the
# worksheet's table relationship (rId1).
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4d0723c7b4f65786e45ea3ecaeede773fa678e9a/e2e/ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.md", "", "", "ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4d0723c7b4f65786e45ea3ecaeede773fa678e9a/.localization-config", "", "", ".localization-config") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4d0723c7b4f65786e45ea3ecaeede773fa678e9a/e2e/ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.md", "", "", "ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/af4d64ec1fc8ee1f13c1d060947f617bb5b087df/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.a40d413bfd1eab2c18f776a1e19828eda35e95b9.zh-cn.xlf", "", "", "ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.a40d413bfd1eab2c18f776a1e19828eda35e95b9.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/cc87e14258571e690db22b45c7d4b00c01cb2308/e2e/ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.md", "", "", "ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/06db303502902256b034ec73fc2cea85564b41c3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.a40d413bfd1eab2c18f776a1e19828eda35e95b9.zh-cn.xlf", "", "", "ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.a40d413bfd1eab2c18f776a1e19828eda35e95b9.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4d0723c7b4f65786e45ea3ecaeede773fa678e9a/.localization-config", "", "", ".localization-config") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4d0723c7b4f65786e45ea3ecaeede773fa678e9a/e2e/ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.md", "", "", "ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a31719d1f5ce77212b6e59951750c1e081aa6edb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.a40d413bfd1eab2c18f776a1e19828eda35e95b9.de-de.xlf", "", "", "ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.a40d413bfd1eab2c18f776a1e19828eda35e95b9.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5ee113d96ab0e5b1068904aee861c6ce482b9442/e2e/ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.md", "", "", "ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ecf7bfb7ac4273f7e45b0f98a58043db70525959/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.a40d413bfd1eab2c18f776a1e19828eda35e95b9.de-de.xlf", "", "", "ebee5ea8-fce3-4e2e-8c60-50f1dae031a9.a40d413bfd1eab2c18f776a1e19828eda35e95b9.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4d0723c7b4f65786e45ea3ecaeede773fa678e9a/.localization-config", "", "", ".localization-config") | Out-Null
